# Reposition/resize the five AWS-architecture container shapes (and their
# paired corner icons) on slide 1, per the target OOXML diff.
#
# Note: Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
# and the host truncates (floors) pts*12700 back to EMU, so the literals below
# are chosen as the midpoint of the EMU bucket that floors to the exact target
# EMU value (avoids the naive emu/12700 landing 1 EMU short after truncation).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 5 (id 35, "Rectangle 34" / "AWS Cloud" outer box)
# off 737056,527265 ext 9728104,7048285 -> off 812800,403226 ext 9652360,7172325
$shpAwsCloud = $s.Shapes.Item(5)
$shpAwsCloud.Left   = 64.00003937007874
$shpAwsCloud.Top    = 31.750118110236222
$shpAwsCloud.Width  = 760.0283858267717
$shpAwsCloud.Height = 564.7500393700788

# Shape 6 (id 12, "Graphic 11", AWS Cloud corner icon)
# off 737057,530130 (ext unchanged) -> off 814400,403227
$shpAwsCloudIcon = $s.Shapes.Item(6)
$shpAwsCloudIcon.Left = 64.12602362204724
$shpAwsCloudIcon.Top  = 31.7501968503937

# Shape 7 (id 38, "Rectangle 37" / "Availability Zone 1" box)
# off 1895678,1080117 ext 2680838,6279686 -> off 1895678,711205 ext 2680838,6648598
$shpAz1 = $s.Shapes.Item(7)
$shpAz1.Top    = 56.00043307086614
$shpAz1.Height = 523.5116929133858

# Shape 8 (id 42, "Rectangle 41" / "VPC" box)
# off 954676,965284 ext 8010164,6494882 -> off 1456076,1071096 ext 7508763,6389069
$shpVpc = $s.Shapes.Item(8)
$shpVpc.Left   = 114.65169291338583
$shpVpc.Top    = 84.33830708661418
$shpVpc.Width  = 591.241220472441
$shpVpc.Height = 503.0763385826772

# Shape 9 (id 45, "Graphic 44", VPC corner icon)
# off 957607,967036 (ext unchanged) -> off 1456076,1071096
$shpVpcIcon = $s.Shapes.Item(9)
$shpVpcIcon.Left = 114.65169291338583
$shpVpcIcon.Top  = 84.33830708661418

# Shape 12 (id 67, "Rectangle 66" / "Availability Zone 2" box)
# off 6152819,1080117 ext 2662926,6279685 -> off 6152819,711205 ext 2662926,6648598
$shpAz2 = $s.Shapes.Item(12)
$shpAz2.Top    = 56.00043307086614
$shpAz2.Height = 523.5116929133858
